$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 141; EnName = "Al Oyoun";  ArName = "العيون";                Lat = 25.674498;          Lon = 49.551423;          Area = "منطقة مكة المكرمة";     Region = "غرب المملكة" },
    @{ Row = 142; EnName = "Thowal";    ArName = "ثول";                   Lat = 22.276098999999999; Lon = 39.112707999999998; Area = "المنطقة الشرقية";       Region = "شرق المملكة" },
    @{ Row = 143; EnName = "Turaf";     ArName = "طريف";                  Lat = 31.665737;          Lon = 38.661921999999997; Area = "منطقة الحدود الشمالية"; Region = "شمال المملكة" },
    @{ Row = 144; EnName = "Al Rafaia"; ArName = "الرفيعة (رفائع الجمش)"; Lat = 24.679561;          Lon = 43.684036999999996; Area = "منطقة الرياض";          Region = "وسط المملكة" }
)

# Row 140 is the last pre-existing data row; its formatting (border style) is
# what the new rows should inherit.
$formatSourceRow = 140

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.EnName
    $ws.Cells.Item($row, 2).Value = $r.EnName
    $ws.Cells.Item($row, 3).Value = $r.ArName
    $ws.Cells.Item($row, 4).Value = $r.Lat
    $ws.Cells.Item($row, 5).Value = $r.Lon
    $ws.Cells.Item($row, 6).Value = $r.Area
    $ws.Cells.Item($row, 7).Value = $r.Region

    $src = $ws.Range("A" + $formatSourceRow + ":G" + $formatSourceRow)
    $dst = $ws.Range("A" + $row + ":G" + $row)
    $src.Copy()
    $dst.PasteSpecial(-4122)
}

$ws.Application.CutCopyMode = $false
